# Actualización desde MV -datos-
# Adds a new "Agosto.2021" vintage column (BH) to the quarterly series table,
# carrying forward the last published estimate for each existing row, revises
# the most recent quarter's estimate, and appends a brand-new row for the
# 01-04-2021 quarter which is only observed in this newest vintage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell (BH1), with the same look & feel as the rest of row 1 ---
$ws.Range("BH1").Value = "Agosto.2021"
$ws.Range("BG1").Copy() | Out-Null
$ws.Range("BH1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Carry the last available estimate (column BG) forward into the new
#     column BH for every existing data row except the last one ---
for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 59).Copy() | Out-Null        # BG = column 59
    $ws.Cells.Item($r, 60).PasteSpecial(-4163) | Out-Null  # xlPasteValues -> BH
}

# --- Row 74 (quarter 01-01-2021) gets a revised estimate in the new column ---
$ws.Range("BH74").Value = 25311

# --- Brand-new row 75 for the 01-04-2021 quarter, first observed this vintage ---
# Force text so Excel's smart-entry parser doesn't turn the dd-mm-yyyy-looking
# label into a real date (matches the "Serie" labels used by every other row),
# then drop the format back to Normal so the cell keeps the sheet's default look.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"
$ws.Range("BH75").Value = 26556
